$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells whose new value looks like a plain number need to be
# forced to remain Text (matching the original inlineStr/"General"-looking
# text cells) - otherwise COM auto-converts them to real numbers and we lose
# the exact decimal text (e.g. "249.98" -> 249.97999999999999, "0.120" -> 0.12).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.548.15"
$ws.Range("E2").Value = "  -1.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.221.24"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.30%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "249.98"
$ws.Range("E5").Value = "  +7.46%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.625"
$ws.Range("E6").Value = "  -1.08%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "70.66"
$ws.Range("E7").Value = "  +1.34%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.08%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.567"
$ws.Range("E9").Value = "  +2.42%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "41.67"
$ws.Range("E10").Value = "  +16.58%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0956"
$ws.Range("E11").Value = "  -2.98%  "

# Row 12 - OKB
Set-TextValue $ws.Range("D12") "58.65"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.02%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.98"
$ws.Range("E14").Value = "  +3.05%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.555.37"
$ws.Range("E15").Value = "  -1.06%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "14.84"
$ws.Range("E16").Value = "  -1.40%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.852"
$ws.Range("E17").Value = "  -1.11%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.220.69"
$ws.Range("E18").Value = "  -1.07%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "41.468.49"
$ws.Range("E19").Value = "  -1.10%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0962"
$ws.Range("E20").Value = "  -0.97%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.19"
$ws.Range("E21").Value = "  -0.80%  "

# Row 22 - Litecoin
Set-TextValue $ws.Range("D22") "72.59"
$ws.Range("E22").Value = "  -1.05%  "

# Row 23 - ImmutableX
Set-TextValue $ws.Range("D23") "2.24"
$ws.Range("E23").Value = "  +9.90%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "233.87"
$ws.Range("E24").Value = "  -1.29%  "

# Row 25 - WEMIXToken
Set-TextValue $ws.Range("D25") "3.83"
$ws.Range("E25").Value = "  +4.97%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.03%  "

# Row 27 - PancakeSwap
Set-TextValue $ws.Range("D27") "2.49"
$ws.Range("E27").Value = "  +5.51%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("D28") "10.58"
$ws.Range("E28").Value = "  +5.57%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +0.33%  "

# Row 30 - Monero
Set-TextValue $ws.Range("D30") "170.97"
$ws.Range("E30").Value = "  +0.98%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Range("D31") "20.56"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32 - Kaspa
Set-TextValue $ws.Range("D32") "0.120"
$ws.Range("E32").Value = "  +0.70%  "

# Row 33 - Stellar
Set-TextValue $ws.Range("D33") "0.124"
$ws.Range("E33").Value = "  -2.69%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D34") "5.53"
$ws.Range("E34").Value = "  +0.51%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0718"
$ws.Range("E35").Value = "  +0.34%  "

# Row 36 - was InjectiveProtocol, now Filecoin (rows 36/37 swapped content)
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D36") "4.67"
$ws.Range("E36").Value = "  -2.11%  "

# Row 37 - was Filecoin, now InjectiveProtocol
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D37") "26.07"
$ws.Range("E37").Value = "  +19.04%  "

# Row 38 - RenderToken
Set-TextValue $ws.Range("D38") "3.94"
$ws.Range("E38").Value = "  +9.43%  "

# Row 39 - VeChain
Set-TextValue $ws.Range("D39") "0.0290"
$ws.Range("E39").Value = "  +9.08%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  +1.29%  "

# Row 41 - MultiversX
Set-TextValue $ws.Range("D41") "68.80"
$ws.Range("E41").Value = "  +3.13%  "

# Row 42 - THORChain
Set-TextValue $ws.Range("D42") "5.94"
$ws.Range("E42").Value = "  -1.57%  "

# Row 43 - Celestia
Set-TextValue $ws.Range("D43") "11.99"
$ws.Range("E43").Value = "  +21.35%  "

# Row 44 - Algorand
Set-TextValue $ws.Range("D44") "0.207"
$ws.Range("E44").Value = "  +8.34%  "

# Row 45 - FTXToken
Set-TextValue $ws.Range("D45") "4.88"
$ws.Range("E45").Value = "  -2.79%  "

# Row 46 - was SynthetixNetwork, now FraxShare (rows 46/47 swapped content)
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D46") "8.72"
$ws.Range("E46").Value = "  -4.12%  "

# Row 47 - was FraxShare, now SynthetixNetwork
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D47") "4.74"
$ws.Range("E47").Value = "  +8.60%  "

# Row 48 - Cronos
Set-TextValue $ws.Range("D48") "0.101"
$ws.Range("E48").Value = "  +0.28%  "

# Row 49 - BinanceUSD
$ws.Range("E49").Value = "  -0.30%  "

# Row 50 - ARBITRUM
$ws.Range("E50").Value = "  +6.87%  "

# Row 51 - TrustWalletToken
Set-TextValue $ws.Range("D51") "1.19"
$ws.Range("E51").Value = "  +1.41%  "
